$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that failed ("Fallo" -> profit -1)
$falloRows = @(83, 92, 95, 96, 115, 117)
foreach ($r in $falloRows) {
    $ws.Range("G$r").Value = "Fallo"
    $ws.Range("H$r").Value = -1
}

# Row 107 succeeded ("Acierto" -> profit equals the cuota, 2.25)
$ws.Range("G107").Value = "Acierto"
$ws.Range("H107").Value = 2.25

# Rows 121-124: event_id (column A) should be stored as a number, not text
$ws.Range("A121").Value = 14552518
$ws.Range("A122").Value = 14552529
$ws.Range("A123").Value = 14552909
$ws.Range("A124").Value = 14552660
